$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table of cell updates (row/col level deltas) derived from the
# refreshed cryptocurrency market data snapshot.
$updates = @(
    @{ Cell = 'D2'; Value = '30.323.21' },
    @{ Cell = 'E2'; Value = '  +0.09%  ' },
    @{ Cell = 'D3'; Value = '1.929.44' },
    @{ Cell = 'D4'; Value = '1.003' },
    @{ Cell = 'E4'; Value = '  +0.29%  ' },
    @{ Cell = 'D5'; Value = '0.7448' },
    @{ Cell = 'E5'; Value = '  +3.89%  ' },
    @{ Cell = 'D6'; Value = '244.27' },
    @{ Cell = 'E6'; Value = '  -1.93%  ' },
    @{ Cell = 'D7'; Value = '1.003' },
    @{ Cell = 'E7'; Value = '  +0.30%  ' },
    @{ Cell = 'D8'; Value = '0.3157' },
    @{ Cell = 'E8'; Value = '  -1.46%  ' },
    @{ Cell = 'D9'; Value = '27.43' },
    @{ Cell = 'E9'; Value = '  -2.24%  ' },
    @{ Cell = 'D10'; Value = '0.07059' },
    @{ Cell = 'E10'; Value = '  -0.53%  ' },
    @{ Cell = 'D11'; Value = '0.7792' },
    @{ Cell = 'E11'; Value = '  -1.37%  ' },
    @{ Cell = 'D12'; Value = '0.08048' },
    @{ Cell = 'E12'; Value = '  +0.60%  ' },
    @{ Cell = 'B13'; Value = 'Polkadot' },
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' },
    @{ Cell = 'D13'; Value = '5.378' },
    @{ Cell = 'E13'; Value = '  +0.01%  ' },
    @{ Cell = 'B14'; Value = 'WrappedEther' },
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' },
    @{ Cell = 'D14'; Value = '1.857.78' },
    @{ Cell = 'E14'; Value = '  -3.66%  ' },
    @{ Cell = 'D15'; Value = '93.17' },
    @{ Cell = 'E15'; Value = '  -1.71%  ' },
    @{ Cell = 'D16'; Value = '14.54' },
    @{ Cell = 'E16'; Value = '  -0.77%  ' },
    @{ Cell = 'D17'; Value = '30.346.24' },
    @{ Cell = 'E17'; Value = '  +0.20%  ' },
    @{ Cell = 'D18'; Value = '5.997' },
    @{ Cell = 'E18'; Value = '  +4.03%  ' },
    @{ Cell = 'D19'; Value = '251.23' },
    @{ Cell = 'E19'; Value = '  -2.29%  ' },
    @{ Cell = 'D20'; Value = '0.000007907' },
    @{ Cell = 'E20'; Value = '  -2.25%  ' },
    @{ Cell = 'D21'; Value = '2.182.75' },
    @{ Cell = 'E21'; Value = '  +0.00%  ' },
    @{ Cell = 'D22'; Value = '1.002' },
    @{ Cell = 'E22'; Value = '  +0.18%  ' },
    @{ Cell = 'E23'; Value = '  +0.17%  ' },
    @{ Cell = 'D24'; Value = '6.641' },
    @{ Cell = 'E24'; Value = '  -2.74%  ' },
    @{ Cell = 'D25'; Value = '9.564' },
    @{ Cell = 'E25'; Value = '  +0.24%  ' },
    @{ Cell = 'D26'; Value = '165.37' },
    @{ Cell = 'E26'; Value = '  +0.49%  ' },
    @{ Cell = 'D27'; Value = '19.03' },
    @{ Cell = 'E27'; Value = '  -0.40%  ' },
    @{ Cell = 'D28'; Value = '0.1289' },
    @{ Cell = 'E28'; Value = '  +1.14%  ' },
    @{ Cell = 'D29'; Value = '2.171' },
    @{ Cell = 'E29'; Value = '  -4.79%  ' },
    @{ Cell = 'D30'; Value = '1.575' },
    @{ Cell = 'E30'; Value = '  +2.88%  ' },
    @{ Cell = 'D31'; Value = '1.358' },
    @{ Cell = 'E31'; Value = '  +0.34%  ' },
    @{ Cell = 'D32'; Value = '4.404' },
    @{ Cell = 'E32'; Value = '  -0.15%  ' },
    @{ Cell = 'D33'; Value = '4.127' },
    @{ Cell = 'E33'; Value = '  -0.35%  ' },
    @{ Cell = 'D34'; Value = '0.05217' },
    @{ Cell = 'E34'; Value = '  +1.79%  ' },
    @{ Cell = 'D35'; Value = '1.310' },
    @{ Cell = 'E35'; Value = '  +2.78%  ' },
    @{ Cell = 'D36'; Value = '0.7535' },
    @{ Cell = 'E36'; Value = '  +0.92%  ' },
    @{ Cell = 'D37'; Value = '2.767' },
    @{ Cell = 'E37'; Value = '  -0.32%  ' },
    @{ Cell = 'D38'; Value = '0.01953' },
    @{ Cell = 'E38'; Value = '  -1.32%  ' },
    @{ Cell = 'D39'; Value = '2.797' },
    @{ Cell = 'E39'; Value = '  -0.08%  ' },
    @{ Cell = 'D40'; Value = '6.507' },
    @{ Cell = 'E40'; Value = '  +2.05%  ' },
    @{ Cell = 'D41'; Value = '76.70' },
    @{ Cell = 'E41'; Value = '  -1.78%  ' },
    @{ Cell = 'D42'; Value = '0.4505' },
    @{ Cell = 'E42'; Value = '  -0.10%  ' },
    @{ Cell = 'D43'; Value = '1.964' },
    @{ Cell = 'E43'; Value = '  -1.53%  ' },
    @{ Cell = 'D44'; Value = '0.8429' },
    @{ Cell = 'E44'; Value = '  -0.23%  ' },
    @{ Cell = 'D45'; Value = '1.002' },
    @{ Cell = 'E45'; Value = '  +0.21%  ' },
    @{ Cell = 'D46'; Value = '10.01' },
    @{ Cell = 'E46'; Value = '  +2.02%  ' },
    @{ Cell = 'D47'; Value = '7.663' },
    @{ Cell = 'E47'; Value = '  +2.48%  ' },
    @{ Cell = 'D48'; Value = '101.43' },
    @{ Cell = 'E48'; Value = '  +0.50%  ' },
    @{ Cell = 'B49'; Value = 'Elrond' },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld' },
    @{ Cell = 'D49'; Value = '37.70' },
    @{ Cell = 'E49'; Value = '  +2.34%  ' },
    @{ Cell = 'B50'; Value = 'RocketPoolETH' },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth' },
    @{ Cell = 'D50'; Value = '2.060.56' },
    @{ Cell = 'E50'; Value = '  -1.23%  ' },
    @{ Cell = 'D51'; Value = '0.1226' },
    @{ Cell = 'E51'; Value = '  +7.20%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    # Force text storage so numeric-looking strings (e.g. "30.323.21",
    # "76.70") are preserved verbatim instead of being coerced into
    # Excel numbers, then restore the default "Normal" style so no
    # stray quote-prefix / number-format markers are left behind.
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
